# Auto-generated edit script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.442.11"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "2.098.68"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.26"
$ws.Range("E5").Value = "  +1.50%  "
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5208"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4523"
$ws.Range("E8").Value = "  +3.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.57"
$ws.Range("E9").Value = "  +15.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08870"
$ws.Range("E10").Value = "  +0.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.178"
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.99"
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("D13").Value = "2.102.12"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.791"
$ws.Range("E14").Value = "  +0.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.991"
$ws.Range("E15").Value = "  +2.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.70"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.007"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001140"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06636"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.13"
$ws.Range("E20").Value = "  +0.64%  "
$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.286"
$ws.Range("E22").Value = "  -0.96%  "
$ws.Range("D23").Value = "30.505.22"
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.35"
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.338"
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("D26").Value = "2.348.33"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.11"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.47"
$ws.Range("E28").Value = "  +0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.512"
$ws.Range("E29").Value = "  -2.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.12"
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.198"
$ws.Range("E31").Value = "  -1.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1065"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.646"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.397"
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.949"
$ws.Range("E35").Value = "  +0.74%  "
$ws.Range("E36").Value = "  +4.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.769"
$ws.Range("E37").Value = "  +5.02%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02572"
$ws.Range("E38").Value = "  -0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06829"
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2293"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.67"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6840"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.243"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.313"
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.91"
$ws.Range("E45").Value = "  -0.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6334"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.666"
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.245"
$ws.Range("E48").Value = "  -0.69%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000347"
$ws.Range("E49").Value = "  +17.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.202"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "82.91"
$ws.Range("E51").Value = "  +0.36%  "
